$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lowercase all the header values in row 1 (A1:J1), keeping column order the same
$ws.Range("A1").Value = "sample_id"
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "lifestage"
$ws.Range("D1").Value = "photo"
$ws.Range("E1").Value = "diameter"
$ws.Range("F1").Value = "date_collected"
$ws.Range("G1").Value = "site"
$ws.Range("H1").Value = "tube_id"
$ws.Range("I1").Value = "experiment"
$ws.Range("J1").Value = "notes"

# Update the active selection to B2
$ws.Range("B2").Select()
